$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2008.4
$ws.Range("J28").Value = 2492
$ws.Range("L28").Value = 2492
$ws.Range("N28").Value = -3462
$ws.Range("H63").Value = 230000
$ws.Range("J63").Value = 230000
$ws.Range("L63").Value = 230000
$ws.Range("N63").Value = -231248
$ws.Range("H66").Value = 230000
$ws.Range("J66").Value = 230000
$ws.Range("L66").Value = 690000
$ws.Range("N66").Value = -696240
$ws.Range("H80").Value = 58831680
$ws.Range("I80").Value = 200000960
$ws.Range("J80").Value = 11147.667
$ws.Range("K80").Value = 600002880
$ws.Range("L80").Value = 33443.001
$ws.Range("M80").Value = -600001882
$ws.Range("N80").Value = -35439.001
$ws.Range("H83").Value = 58831680
$ws.Range("I83").Value = 200000960
$ws.Range("J83").Value = 11147.667
$ws.Range("K83").Value = 1800008640
$ws.Range("L83").Value = 100329.003
$ws.Range("M83").Value = -1800003648
$ws.Range("N83").Value = -110313.003
$ws.Range("H98").Value = 1692.95
$ws.Range("I98").Value = 1489.6842
$ws.Range("K98").Value = 1489.6842
$ws.Range("M98").Value = 8.315800000000081
$ws.Range("H107").Value = 1722
$ws.Range("I107").Value = 1948.8636
$ws.Range("J107").Value = 474.25
$ws.Range("K107").Value = 1948.8636
$ws.Range("L107").Value = 474.25
$ws.Range("M107").Value = -28.86359999999991
$ws.Range("N107").Value = -4314.25
$ws.Range("H116").Value = 83340280
$ws.Range("I116").Value = 125005960
$ws.Range("K116").Value = 125005960
$ws.Range("M116").Value = -125002518
$ws.Range("H122").Value = 1692.95
$ws.Range("I122").Value = 1489.6842
$ws.Range("K122").Value = 4469.0526
$ws.Range("M122").Value = -2019.0526
$ws.Range("H137").Value = 19225.428
$ws.Range("I137").Value = 23896.184
$ws.Range("J137").Value = 3090.0908
$ws.Range("K137").Value = 71688.552
$ws.Range("L137").Value = 9270.2724
$ws.Range("M137").Value = -69138.552
$ws.Range("N137").Value = -14370.2724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4538.636
$ws.Range("I45").Value = 4653.125
$ws.Range("J45").Value = 4233.3335
$ws.Range("K45").Value = 4653.125
$ws.Range("L45").Value = 4233.3335
$ws.Range("M45").Value = -4276.125
$ws.Range("N45").Value = -4987.3335
$ws.Range("H61").Value = 1717.898
$ws.Range("I61").Value = 1597.4375
$ws.Range("K61").Value = 1597.4375
$ws.Range("M61").Value = -1385.4375
$ws.Range("H74").Value = 2086.6086
$ws.Range("I74").Value = 1975.8334
$ws.Range("K74").Value = 1975.8334
$ws.Range("M74").Value = -1101.8334
$ws.Range("H77").Value = 2086.6086
$ws.Range("I77").Value = 1975.8334
$ws.Range("K77").Value = 9879.166999999999
$ws.Range("M77").Value = -5511.166999999999
$ws.Range("H110").Value = 2273.6667
$ws.Range("I110").Value = 1967.6666
$ws.Range("K110").Value = 1967.6666
$ws.Range("M110").Value = 77.33339999999998
$ws.Range("H132").Value = 151642.47
$ws.Range("I132").Value = 183718.03
$ws.Range("K132").Value = 551154.09
$ws.Range("M132").Value = -548624.09
$ws.Range("H136").Value = 1717.898
$ws.Range("I136").Value = 1597.4375
$ws.Range("K136").Value = 4792.3125
$ws.Range("M136").Value = -2242.3125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1158.6129
$ws.Range("I94").Value = 989.4211
$ws.Range("J94").Value = 1426.5
$ws.Range("K94").Value = 989.4211
$ws.Range("L94").Value = 1426.5
$ws.Range("M94").Value = -538.4211
$ws.Range("N94").Value = -2328.5
$ws.Range("H134").Value = 2166853.8
$ws.Range("I134").Value = 2383022.8
$ws.Range("K134").Value = 7149068.399999999
$ws.Range("M134").Value = -7146533.399999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3102.2058
$ws.Range("I31").Value = 2195.5925
$ws.Range("K31").Value = 2195.5925
$ws.Range("M31").Value = -1900.5925
$ws.Range("H34").Value = 3102.2058
$ws.Range("I34").Value = 2195.5925
$ws.Range("K34").Value = 2195.5925
$ws.Range("M34").Value = -1993.5925
$ws.Range("H58").Value = 2362.348
$ws.Range("I58").Value = 2111.1667
$ws.Range("K58").Value = 2111.1667
$ws.Range("M58").Value = -1908.1667
$ws.Range("H134").Value = 2572.3462
$ws.Range("I134").Value = 2555.4688
$ws.Range("J134").Value = 2599.35
$ws.Range("K134").Value = 7666.4064
$ws.Range("L134").Value = 7798.049999999999
$ws.Range("M134").Value = -5131.4064
$ws.Range("N134").Value = -12868.05
$ws.Range("H136").Value = 2362.348
$ws.Range("I136").Value = 2111.1667
$ws.Range("K136").Value = 6333.500100000001
$ws.Range("M136").Value = -3783.500100000001
$ws.Range("H141").Value = 249359.47
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 249359.47
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 249359.47
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -259719.47

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 814.0714
$ws.Range("J44").Value = 932.5
$ws.Range("L44").Value = 2797.5
$ws.Range("N44").Value = -3593.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4445
$ws.Range("I70").Value = 4128.75
$ws.Range("K70").Value = 4128.75
$ws.Range("M70").Value = -3858.75
$ws.Range("H73").Value = 4445
$ws.Range("I73").Value = 4128.75
$ws.Range("K73").Value = 4128.75
$ws.Range("M73").Value = -3192.75
$ws.Range("H107").Value = 860.4
$ws.Range("I107").Value = 940.55554
$ws.Range("K107").Value = 940.55554
$ws.Range("M107").Value = 979.44446
$ws.Range("H113").Value = 1972
$ws.Range("I113").Value = 1314.8334
$ws.Range("J113").Value = 2629.1667
$ws.Range("K113").Value = 1314.8334
$ws.Range("L113").Value = 2629.1667
$ws.Range("M113").Value = 855.1666
$ws.Range("N113").Value = -6969.1667
$ws.Range("H132").Value = 3873.568
$ws.Range("I132").Value = 4018.3928
$ws.Range("K132").Value = 12055.1784
$ws.Range("M132").Value = -9525.178400000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6710.2607
$ws.Range("I7").Value = 5835.8
$ws.Range("K7").Value = 5835.8
$ws.Range("M7").Value = -5723.8
$ws.Range("H43").Value = 1000000000
$ws.Range("J43").Value = 1000000000
$ws.Range("L43").Value = 1000000000
$ws.Range("N43").Value = -1000000386
$ws.Range("H46").Value = 5828.1284
$ws.Range("J46").Value = 4289.3335
$ws.Range("L46").Value = 4289.3335
$ws.Range("N46").Value = -4665.3335
$ws.Range("H61").Value = 2661.4866
$ws.Range("I61").Value = 1653.4828
$ws.Range("K61").Value = 1653.4828
$ws.Range("M61").Value = -1451.4828
$ws.Range("H113").Value = 2661.4866
$ws.Range("I113").Value = 1653.4828
$ws.Range("K113").Value = 1653.4828
$ws.Range("M113").Value = 516.5172
$ws.Range("H122").Value = 15998.474
$ws.Range("I122").Value = 19716
$ws.Range("J122").Value = 11867.889
$ws.Range("K122").Value = 59148
$ws.Range("L122").Value = 35603.667
$ws.Range("M122").Value = -56698
$ws.Range("N122").Value = -40503.667
$ws.Range("H126").Value = 6710.2607
$ws.Range("I126").Value = 5835.8
$ws.Range("K126").Value = 17507.4
$ws.Range("M126").Value = -15037.4
$ws.Range("H132").Value = 436801.3
$ws.Range("I132").Value = 484594.06
$ws.Range("K132").Value = 1453782.18
$ws.Range("M132").Value = -1451252.18
$ws.Range("H136").Value = 2341.9033
$ws.Range("I136").Value = 1722
$ws.Range("K136").Value = 5166
$ws.Range("M136").Value = -2616

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5187.75
$ws.Range("I62").Value = 3249.6667
$ws.Range("J62").Value = 6350.6
$ws.Range("K62").Value = 3249.6667
$ws.Range("L62").Value = 6350.6
$ws.Range("M62").Value = -2625.6667
$ws.Range("N62").Value = -7598.6
$ws.Range("H65").Value = 5187.75
$ws.Range("I65").Value = 3249.6667
$ws.Range("J65").Value = 6350.6
$ws.Range("K65").Value = 16248.3335
$ws.Range("L65").Value = 31753
$ws.Range("M65").Value = -13128.3335
$ws.Range("N65").Value = -37993
$ws.Range("H96").Value = 9760.357
$ws.Range("I96").Value = 6910.5557
$ws.Range("K96").Value = 6910.5557
$ws.Range("M96").Value = -5537.5557
$ws.Range("H113").Value = 1114.1333
$ws.Range("I113").Value = 1120.125
$ws.Range("J113").Value = 1107.2858
$ws.Range("K113").Value = 3360.375
$ws.Range("L113").Value = 3321.8574
$ws.Range("M113").Value = -1190.375
$ws.Range("N113").Value = -7661.857400000001
$ws.Range("H132").Value = 25368.955
$ws.Range("I132").Value = 32658.969
$ws.Range("K132").Value = 97976.90700000001
$ws.Range("M132").Value = -95446.90700000001
$ws.Range("H136").Value = 18233.334
$ws.Range("I136").Value = 1045.9592
$ws.Range("K136").Value = 3137.8776
$ws.Range("M136").Value = -587.8775999999998
